$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect columns that Excel would otherwise auto-convert to numbers/dates
# (currency like "$69.99" and bare dates like "2024-09-22") by forcing a
# Text number format before writing the literal strings, then clearing the
# explicit formatting again so the new rows keep the same (default) style as
# the existing data rows.
$newRange = $ws.Range("A40:F44")
$newRange.NumberFormat = "@"

$ws.Cells.Item(40, 1).Value = '2024-09-22 23:35:51'
$ws.Cells.Item(40, 2).Value = 'monitor_price'
$ws.Cells.Item(40, 3).Value = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'
$ws.Cells.Item(40, 4).Value = '$69.99'
$ws.Cells.Item(40, 5).Value = '2024-09-22'
$ws.Cells.Item(40, 6).Value = '23:35:51'

$ws.Cells.Item(41, 1).Value = '2024-09-22 23:41:37'
$ws.Cells.Item(41, 2).Value = 'monitor_price'
$ws.Cells.Item(41, 3).Value = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'
$ws.Cells.Item(41, 4).Value = '$69.99'
$ws.Cells.Item(41, 5).Value = '2024-09-22'
$ws.Cells.Item(41, 6).Value = '23:41:37'

$ws.Cells.Item(42, 1).Value = '2024-09-22 23:43:01'
$ws.Cells.Item(42, 2).Value = 'monitor_price'
$ws.Cells.Item(42, 3).Value = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'
$ws.Cells.Item(42, 4).Value = '$69.99'
$ws.Cells.Item(42, 5).Value = '2024-09-22'
$ws.Cells.Item(42, 6).Value = '23:43:01'

$ws.Cells.Item(43, 1).Value = '2024-09-22 23:43:27'
$ws.Cells.Item(43, 2).Value = 'monitor_price'
$ws.Cells.Item(43, 3).Value = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'
$ws.Cells.Item(43, 4).Value = '$69.99'
$ws.Cells.Item(43, 5).Value = '2024-09-22'
$ws.Cells.Item(43, 6).Value = '23:43:27'

$ws.Cells.Item(44, 1).Value = '2024-09-22 23:43:50'
$ws.Cells.Item(44, 2).Value = 'monitor_price'
$ws.Cells.Item(44, 3).Value = 'https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960'
$ws.Cells.Item(44, 4).Value = '$69.99'
$ws.Cells.Item(44, 5).Value = '2024-09-22'
$ws.Cells.Item(44, 6).Value = '23:43:50'

$newRange.ClearFormats()
